$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage for values that would
# otherwise be auto-converted by Excel into a number (e.g. "0.656" -> 0.656).
# The NumberFormat/Style dance avoids leaving a permanent style override on the
# cell (Style = "Normal" resets it back to the default, unstyled cell).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "37.036.26"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "2.042.85"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "245.87"
$ws.Range("E5").Value = "  -1.72%  "
Set-TextValue "D6" "0.656"
$ws.Range("E6").Value = "  -1.86%  "
Set-TextValue "D7" "58.80"
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("E8").Value = "  -0.05%  "
Set-TextValue "D9" "0.377"
$ws.Range("E9").Value = "  -2.10%  "
Set-TextValue "D10" "0.0772"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  +2.12%  "
Set-TextValue "D12" "15.40"
$ws.Range("E12").Value = "  -4.97%  "
Set-TextValue "D13" "0.898"
$ws.Range("E13").Value = "  +9.00%  "
$ws.Range("D14").Value = "2.337.08"
$ws.Range("E14").Value = "  -0.79%  "
Set-TextValue "D15" "5.70"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "2.028.76"
$ws.Range("E16").Value = "  -1.27%  "
Set-TextValue "D17" "18.48"
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").Value = "36.960.84"
$ws.Range("E18").Value = "  -0.55%  "
Set-TextValue "D19" "73.63"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -2.04%  "
Set-TextValue "D21" "5.39"
$ws.Range("E21").Value = "  -0.74%  "
Set-TextValue "D22" "239.36"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +1.74%  "
Set-TextValue "D25" "9.64"
$ws.Range("E25").Value = "  +2.70%  "
Set-TextValue "D26" "168.41"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -3.56%  "
Set-TextValue "D28" "19.97"
$ws.Range("E28").Value = "  -0.23%  "
Set-TextValue "D29" "5.55"
$ws.Range("E29").Value = "  +15.31%  "
$ws.Range("E30").Value = "  -1.02%  "
Set-TextValue "D31" "1.12"
$ws.Range("E31").Value = "  -2.67%  "
Set-TextValue "D32" "4.76"
$ws.Range("E32").Value = "  +4.48%  "
Set-TextValue "D33" "0.0613"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  +6.66%  "
Set-TextValue "D36" "0.0852"
$ws.Range("E36").Value = "  -5.31%  "
Set-TextValue "D37" "2.26"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -4.31%  "
Set-TextValue "D39" "5.23"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D41" "0.0973"
$ws.Range("E41").Value = "  -11.13%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.0222"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D43" "1.15"
$ws.Range("E43").Value = "  +0.55%  "
Set-TextValue "D44" "97.48"
$ws.Range("E44").Value = "  -0.77%  "
Set-TextValue "D45" "16.96"
$ws.Range("E45").Value = "  -7.50%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D46" "2.39"
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.298.24"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "6.73"
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D50" "3.65"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.222.38"
$ws.Range("E51").Value = "  -0.88%  "
